$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the bold/centered header style
# already used by the other header cells (B1:G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats: copy G1's style onto H1

# New "Save" data column values (era data updated).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
